# Sample Project / Main.xlsx - "Rules" sheet
#
# Cell B11 held the rule-name "R40" (a shared string). The edit renames
# that rule to "1": the cell now shows the text "1" instead of "R40".
# Because "1" looks like a number, a plain Range.Value assignment would
# make Excel store it as a numeric literal (and typing it with a leading
# apostrophe would flip on the "number stored as text" quote-prefix cell
# style). Neither of those matches the source edit, which simply swaps in
# a new shared-string text value while leaving the cell's existing
# formatting untouched. Writing the text through a text-formatted helper
# cell and pasting only the value into B11 reproduces that: B11 keeps its
# original style/border and ends up with a genuine text value "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$helper = $ws.Range("ZZ1")
$helper.Value = "'1"
$helper.Copy()

$ws.Cells.Item(11, 2).PasteSpecial(-4163)  # xlPasteValues

$helper.Clear()
$excel.CutCopyMode = $false
